# Fruta / hortaliza, semanal
# Updates the weekly Frambuesa (raspberry) price records for Femacal de La
# Calera: each data row (2-15) gets its date, volume, min/max/weighted
# prices, origin province and $/Kg refreshed to the latest weekly values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44587
$ws.Range("M2").Value = 165
$ws.Range("N2").Value = 6500
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6742
$ws.Range("R2").Value = 'Provincia de Linares'
$ws.Range("S2").Value = 3371

# Row 3
$ws.Range("D3").Value = 44211
$ws.Range("M3").Value = 45
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("S3").Value = 3000

# Row 4
$ws.Range("D4").Value = 44589
$ws.Range("M4").Value = 60
$ws.Range("R4").Value = 'Provincia de Curicó'

# Row 5
$ws.Range("D5").Value = 44585

# Row 6
$ws.Range("D6").Value = 44606
$ws.Range("M6").Value = 45
$ws.Range("R6").Value = 'Provincia de Linares'

# Row 7
$ws.Range("D7").Value = 44586
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 7000
$ws.Range("P7").Value = 7000
$ws.Range("S7").Value = 3500

# Row 8
$ws.Range("D8").Value = 44592
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("R8").Value = 'Provincia de Linares'
$ws.Range("S8").Value = 4000

# Row 9
$ws.Range("D9").Value = 44582
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6500
$ws.Range("P9").Value = 6233
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 3116

# Row 10
$ws.Range("D10").Value = 44588
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 6500
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 6750
$ws.Range("R10").Value = 'Provincia de Curicó'
$ws.Range("S10").Value = 3375

# Row 11
$ws.Range("D11").Value = 44614
$ws.Range("M11").Value = 45
$ws.Range("R11").Value = 'Provincia de Linares'

# Row 12
$ws.Range("D12").Value = 44214
$ws.Range("M12").Value = 48
$ws.Range("N12").Value = 6000
$ws.Range("O12").Value = 6000
$ws.Range("P12").Value = 6000
$ws.Range("S12").Value = 3000

# Row 13
$ws.Range("D13").Value = 44627
$ws.Range("N13").Value = 6000
$ws.Range("O13").Value = 6000
$ws.Range("P13").Value = 6000
$ws.Range("S13").Value = 3000

# Row 14
$ws.Range("D14").Value = 44628
$ws.Range("M14").Value = 40

# Row 15
$ws.Range("D15").Value = 44209
$ws.Range("M15").Value = 58
$ws.Range("R15").Value = 'Provincia de Curicó'
